$wb = $excel.ActiveWorkbook

# --- Sheet1: numeric RF concentration values (column E, rows 2-16) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("E2").Value = 5.21247659109051
$ws1.Range("E3").Value = 716.2248375038637
$ws1.Range("E4").Value = 4356.548997543381
$ws1.Range("E5").Value = 6.951565824544648
$ws1.Range("E6").Value = 6.347851084684235
$ws1.Range("E7").Value = 49.53332119676588
$ws1.Range("E8").Value = 136.2929689729838
$ws1.Range("E9").Value = 32.74871911544427
$ws1.Range("E10").Value = 149.3168866482802
$ws1.Range("E11").Value = 76.96197658911851
$ws1.Range("E12").Value = 7.270352946288932
$ws1.Range("E13").Value = 102.7328934421345
$ws1.Range("E14").Value = 35.52751750611267
$ws1.Range("E15").Value = 452.5153282036359
$ws1.Range("E16").Value = 54.34994527723487

# --- Sheet2: RF concentration text values (column E) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("E2").NumberFormat = "@"
$ws2.Range("E2").Value = "5.21247659109051"
$ws2.Range("E2").Style = "Normal"
$ws2.Range("E6").NumberFormat = "@"
$ws2.Range("E6").Value = "716.2248375038637"
$ws2.Range("E6").Style = "Normal"
$ws2.Range("E10").NumberFormat = "@"
$ws2.Range("E10").Value = "4356.548997543381"
$ws2.Range("E10").Style = "Normal"
$ws2.Range("E14").NumberFormat = "@"
$ws2.Range("E14").Value = "6.951565824544648"
$ws2.Range("E14").Style = "Normal"
$ws2.Range("E18").NumberFormat = "@"
$ws2.Range("E18").Value = "6.347851084684235"
$ws2.Range("E18").Style = "Normal"
$ws2.Range("E22").NumberFormat = "@"
$ws2.Range("E22").Value = "49.53332119676588"
$ws2.Range("E22").Style = "Normal"
$ws2.Range("E26").NumberFormat = "@"
$ws2.Range("E26").Value = "136.29296897298383"
$ws2.Range("E26").Style = "Normal"
$ws2.Range("E30").NumberFormat = "@"
$ws2.Range("E30").Value = "32.748719115444274"
$ws2.Range("E30").Style = "Normal"
$ws2.Range("E34").NumberFormat = "@"
$ws2.Range("E34").Value = "149.3168866482802"
$ws2.Range("E34").Style = "Normal"
$ws2.Range("E38").NumberFormat = "@"
$ws2.Range("E38").Value = "76.9619765891185"
$ws2.Range("E38").Style = "Normal"
$ws2.Range("E42").NumberFormat = "@"
$ws2.Range("E42").Value = "7.270352946288932"
$ws2.Range("E42").Style = "Normal"
$ws2.Range("E46").NumberFormat = "@"
$ws2.Range("E46").Value = "102.73289344213454"
$ws2.Range("E46").Style = "Normal"
$ws2.Range("E50").NumberFormat = "@"
$ws2.Range("E50").Value = "35.527517506112666"
$ws2.Range("E50").Style = "Normal"
$ws2.Range("E54").NumberFormat = "@"
$ws2.Range("E54").Value = "452.5153282036359"
$ws2.Range("E54").Style = "Normal"
$ws2.Range("E58").NumberFormat = "@"
$ws2.Range("E58").Value = "54.349945277234866"
$ws2.Range("E58").Style = "Normal"

# --- Sheet2: Best Match / Similarity text values (column B) ---
$ws2.Range("B4").Value = "Best Match: C(C(F)(F)F)(OC(C(F)(F)S(=O)(=O)O)(F)F)(F)F with Similarity: 0.0458015267175573"
$ws2.Range("B8").Value = "Best Match: C(C(F)(F)F)(OC(C(F)(F)S(=O)(=O)O)(F)F)(F)F with Similarity: 0.0512820512820513"
$ws2.Range("B12").Value = "Best Match: C(C(C(C(F)(F)Cl)(F)F)(F)F)(C(C(OC(C(F)(F)S(=O)(=O)O)(F)F)(F)F)(F)F)(F)F with Similarity: 0.0588235294117647"
$ws2.Range("B16").Value = "Best Match: C(C(F)(F)F)(OC(C(F)(F)S(=O)(=O)O)(F)F)(F)F with Similarity: 0.028169014084507"
$ws2.Range("B20").Value = "Best Match: C(C(F)(F)F)(OC(C(F)(F)S(=O)(=O)O)(F)F)(F)F with Similarity: 0.0222222222222222"
$ws2.Range("B24").Value = "Best Match: C(C(C(C(F)(F)Cl)(F)F)(F)F)(C(C(OC(C(F)(F)S(=O)(=O)O)(F)F)(F)F)(F)F)(F)F with Similarity: 0.0482758620689655"
$ws2.Range("B28").Value = "Best Match: C(C(C(C(F)(F)S(=O)(=O)O)(F)F)(F)F)(C(C(F)(F)F)(F)F)(F)F with Similarity: 0.0188679245283019"
$ws2.Range("B32").Value = "Best Match: C(C(C(C(F)(F)S(=O)(=O)O)(F)F)(F)F)(C(C(F)(F)F)(F)F)(F)F with Similarity: 0.0120481927710843"
$ws2.Range("B36").Value = "Best Match: C(C(C(C(F)(F)Cl)(F)F)(F)F)(C(C(OC(C(F)(F)S(=O)(=O)O)(F)F)(F)F)(F)F)(F)F with Similarity: 0.0506329113924051"
$ws2.Range("B40").Value = "Best Match: C(C(C(C(F)(F)S(=O)(=O)O)(F)F)(F)F)(C(C(F)(F)F)(F)F)(F)F with Similarity: 0.0357142857142857"
$ws2.Range("B44").Value = "Best Match: C(C(F)(F)F)(OC(C(F)(F)S(=O)(=O)O)(F)F)(F)F with Similarity: 0.0428571428571429"
$ws2.Range("B48").Value = "Best Match: C(C(F)(F)F)(OC(C(F)(F)S(=O)(=O)O)(F)F)(F)F with Similarity: 0.0897435897435897"
$ws2.Range("B52").Value = "Best Match: C(C(C(C(F)(F)S(=O)(=O)O)(F)F)(F)F)(C(C(F)(F)F)(F)F)(F)F with Similarity: 0.032967032967033"
$ws2.Range("B60").Value = "Best Match: C(C(C(C(F)(F)S(=O)(=O)O)(F)F)(F)F)(C(C(F)(F)F)(F)F)(F)F with Similarity: 0.0333333333333333"
